# Actualización automática 2025-07-04 17:05:07
#
# Updates July ("julio") sales for client "BORJA TORRES LETTY JANET"
# (advisor "LINDAO ZUÑIGA BRYAN JOSE") in the PORCELANATO category,
# and propagates the new totals through the three report sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": PORCELANATO sales for row 10 (client
# BORJA TORRES LETTY JANET) and the "X de 56" progress counter.
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M10").Value = 1451.52
$wsGrupo.Range("M58").Value = "2 de 56"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL": julio column for the same client/row and
# the julio column total.
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F10").Value = 1451.52
$wsMensual.Range("F58").Value = 2110.51

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": recompute VENTA / POR CUMPLIR /
# CUMPLIMIENTO for the affected categories and the TOTAL row.
# ---------------------------------------------------------------
$wsCump = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 4 - FREGADEROS DE COCINA
$wsCump.Range("D4").Value = 275.59
$wsCump.Range("E4").Value = -133.087904974973
$wsCump.Range("F4").Value = 1.933936479681926

# Row 8 - LAVABOS
$wsCump.Range("D8").Value = 779.4
$wsCump.Range("E8").Value = -29.39999999999998
$wsCump.Range("F8").Value = 1.0392

# Row 16 - PORCELANATO
$wsCump.Range("D16").Value = 2110.51
$wsCump.Range("E16").Value = 38279.66
$wsCump.Range("F16").Value = 0.05225306058380047

# Row 17 - PUERTAS DE SEGURIDAD
$wsCump.Range("D17").Value = 372.66
$wsCump.Range("E17").Value = -30.66000000000003
$wsCump.Range("F17").Value = 1.089649122807018

# Row 19 - TOTAL
$wsCump.Range("D19").Value = 3538.16
$wsCump.Range("E19").Value = 51871.54560036207
$wsCump.Range("F19").Value = 0.06385451721253831

# Column width tweaks (D/E widened slightly to fit new values).
# NOTE: Excel's ColumnWidth (character units) maps to a stored width of
# ColumnWidth + 5/6 in the saved xlsx, so we back that padding out here
# to land on the exact target stored widths of 13 and 24.
$wsCump.Columns.Item(4).ColumnWidth = 13 - 5/6
$wsCump.Columns.Item(5).ColumnWidth = 24 - 5/6

Write-Output "applied"
